# Atualiza base_barcos_dummy.xlsx com novos dados de veleiros.
# Adiciona dois novos barcos (linhas 34 e 35) e ajusta a pagina/selecao
# como ficou registrado apos a edicao no Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Novo barco 33: Capitão Cloroquina --------------------------------
$ws.Range("A34").Value2 = 33
$ws.Range("B34").Value2 = 1234
$ws.Range("C34").Value2 = "Capitão Cloroquina"
$ws.Range("D34").Value2 = "Urca - Rio de Janeiro"
$ws.Range("E34").Value2 = 25
$ws.Range("F34").Value2 = 8
$ws.Range("G34").Value2 = 512
$ws.Range("H34").Value2 = "Churrasco"

# --- Novo barco 34: Amor Lindo -----------------------------------------
$ws.Range("A35").Value2 = 34
$ws.Range("B35").Value2 = 4321
$ws.Range("C35").Value2 = "Amor Lindo"
$ws.Range("D35").Value2 = "Paraty"
$ws.Range("E35").Value2 = 15
$ws.Range("F35").Value2 = 3
$ws.Range("G35").Value2 = 100
$ws.Range("H35").Value2 = "Jantar romântico, DJ a bordo"

# Celula vazia deixada logo abaixo, sublinhada (estado em que o autor
# parou de digitar), e selecionada ao salvar.
$ws.Range("B36").Font.Underline = 1
$ws.Range("B36").Select() | Out-Null

# Ajustes de layout de impressao feitos junto com a edição.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
